$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("223", $false, "13:28"),
    @("737", $false, "13:29"),
    @("712", $false, "13:30"),
    @("652", $false, "13:31"),
    @("675", $false, "13:32"),
    @("507", $false, "13:33"),
    @("639", $false, "13:34"),
    @("1",   $false, "13:35"),
    @("293", $true,  "13:28"),
    @("676", $true,  "13:29"),
    @("716", $true,  "13:30"),
    @("664", $true,  "13:31"),
    @("585", $true,  "13:32"),
    @("627", $true,  "13:33"),
    @("593", $true,  "13:34")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = "'" + $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = "'" + $rec[2]
    $row++
}
